$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9999999914007518
$ws.Range("A2").Value = 0.99475998531723842
$ws.Range("A3").Value = 0.97631595320642628
$ws.Range("A4").Value = 0.9695916255189081
$ws.Range("A5").Value = 0.95987275642605563
$ws.Range("A6").Value = 0.93685751408343576
$ws.Range("A7").Value = 0.93685395429274809
$ws.Range("A8").Value = 0.9349218924990641
$ws.Range("A9").Value = 0.93688932201457686
$ws.Range("A10").Value = 0.93982413047881463
$ws.Range("A11").Value = 0.94019577070061244
$ws.Range("A12").Value = 0.94116162545614968
$ws.Range("A13").Value = 0.95012825561264647
$ws.Range("A14").Value = 0.94861238382070612
$ws.Range("A15").Value = 0.94877363503023693
$ws.Range("A16").Value = 0.94967712874654819
$ws.Range("A17").Value = 0.96074419423054658
$ws.Range("A18").Value = 0.96741044998109249
$ws.Range("A19").Value = 0.99136034897723713
$ws.Range("A20").Value = 0.96336468850113599
$ws.Range("A21").Value = 0.95579621302537632
$ws.Range("A22").Value = 0.95453169604970123
$ws.Range("A23").Value = 0.96782268920439329
$ws.Range("A24").Value = 0.95169971525483288
$ws.Range("A25").Value = 0.94219818322133853
$ws.Range("A26").Value = 0.93320155407963146
$ws.Range("A27").Value = 0.92835434160491004
$ws.Range("A28").Value = 0.9068793312410155
$ws.Range("A29").Value = 0.89160727939589224
$ws.Range("A30").Value = 0.88503678273418274
$ws.Range("A31").Value = 0.87738302234042753
$ws.Range("A32").Value = 0.87570373151831904
$ws.Range("A33").Value = 0.87518373532537241
